$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1676
$ws1.Range("F6").Value = 459
$ws1.Range("F8").Value = 73
$ws1.Range("F9").Value = 589
$ws1.Range("F10").Value = 406

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1676
$ws4.Range("F6").Value = 459
$ws4.Range("F8").Value = 73
$ws4.Range("F9").Value = 589
$ws4.Range("F10").Value = 406
